$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 19:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1111236
$ws.Range("C4").Value = 16213
$ws.Range("D4").Value = 157809
$ws.Range("E4").Value = 888551
$ws.Range("F4").Value = 15118
$ws.Range("G4").Value = 1020
$ws.Range("H4").Value = 64876

# Row 5 - España
$ws.Range("C5").Value = 3648

# Row 9 - Alemania
$ws.Range("E9").Value = 30002
$ws.Range("G9").Value = 17
$ws.Range("H9").Value = 6640

# Row 25 - Irlanda
$ws.Range("B25").Value = 20833
$ws.Range("C25").Value = 221
$ws.Range("E25").Value = 6182
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 1265

# Row 27 - Pakistan
$ws.Range("B27").Value = 17699
$ws.Range("C27").Value = 1226
$ws.Range("E27").Value = 12940
$ws.Range("G27").Value = 47
$ws.Range("H27").Value = 408

# Row 71 - Uzbekistan
$ws.Range("D71").Value = 1212
$ws.Range("E71").Value = 854

# Row 80 - Cuba
$ws.Range("B80").Value = 1537
$ws.Range("C80").Value = 36
$ws.Range("D80").Value = 714
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 64

# Row 123 - Estado de Palestina
$ws.Range("B123").Value = 353
$ws.Range("C123").Value = 9
$ws.Range("E123").Value = 275

# Row 139 - Birmania
$ws.Range("D139").Value = 31
$ws.Range("E139").Value = 114

# Row 161 - Mozambique
$ws.Range("B161").Value = 79
$ws.Range("C161").Value = 3
$ws.Range("E161").Value = 67

$wb.Save()
